# Update column F (dSF) values for specific rows to match the re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = 0
    11 = 3
    16 = 5
    20 = 0
    22 = -6
    23 = 0
    24 = 1
    34 = 4
    35 = 2
    36 = 0
    40 = -2
    44 = -4
    49 = 0
    50 = 1
    59 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
